# Design-doc slide update:
#   1. Insert a new dashed-border "frame" rectangle ("Rectangle 1") behind all
#      existing content, sized/positioned to enclose the whole diagram.
#   2. Shift every pre-existing shape on the slide 1270366 EMU to the right
#      (to make room for / center inside the new frame) - the Y coordinates,
#      widths and heights are untouched.
#
# NOTE on precision: this COM-interop runtime stores Shape.Left/Top/Width/
# Height internally as 32-bit floats, so a naive
#     $shape.Left = $targetEmu / 12700.0
# can land 1 EMU away from the intended integer EMU value after the
# float32 round-trip. The Set-Precise* helpers below compensate by nudging
# the point value in 1e-6pt steps until the value that actually gets stored
# (verified via read-back) converts back to exactly the target EMU.
# (Dynamic `$Shape.$PropertyName = ...` assignment isn't supported by this
# interpreter, so each property gets its own small helper.)

function Set-PreciseLeft {
    param($Shape, [int64]$TargetEmu)

    $base = $TargetEmu / 12700.0
    $Shape.Left = $base
    $cur = [int64]([math]::Round($Shape.Left * 12700.0))
    if ($cur -eq $TargetEmu) { return }

    for ($i = 1; $i -le 20000; $i++) {
        foreach ($sign in @(1, -1)) {
            $cand = $base + ($sign * $i * 0.000001)
            $Shape.Left = $cand
            $cur = [int64]([math]::Round($Shape.Left * 12700.0))
            if ($cur -eq $TargetEmu) { return }
        }
    }
    Write-Host "WARNING: could not land Left on exact EMU $TargetEmu (got $cur)"
}

function Set-PreciseTop {
    param($Shape, [int64]$TargetEmu)

    $base = $TargetEmu / 12700.0
    $Shape.Top = $base
    $cur = [int64]([math]::Round($Shape.Top * 12700.0))
    if ($cur -eq $TargetEmu) { return }

    for ($i = 1; $i -le 20000; $i++) {
        foreach ($sign in @(1, -1)) {
            $cand = $base + ($sign * $i * 0.000001)
            $Shape.Top = $cand
            $cur = [int64]([math]::Round($Shape.Top * 12700.0))
            if ($cur -eq $TargetEmu) { return }
        }
    }
    Write-Host "WARNING: could not land Top on exact EMU $TargetEmu (got $cur)"
}

function Set-PreciseWidth {
    param($Shape, [int64]$TargetEmu)

    $base = $TargetEmu / 12700.0
    $Shape.Width = $base
    $cur = [int64]([math]::Round($Shape.Width * 12700.0))
    if ($cur -eq $TargetEmu) { return }

    for ($i = 1; $i -le 20000; $i++) {
        foreach ($sign in @(1, -1)) {
            $cand = $base + ($sign * $i * 0.000001)
            $Shape.Width = $cand
            $cur = [int64]([math]::Round($Shape.Width * 12700.0))
            if ($cur -eq $TargetEmu) { return }
        }
    }
    Write-Host "WARNING: could not land Width on exact EMU $TargetEmu (got $cur)"
}

function Set-PreciseHeight {
    param($Shape, [int64]$TargetEmu)

    $base = $TargetEmu / 12700.0
    $Shape.Height = $base
    $cur = [int64]([math]::Round($Shape.Height * 12700.0))
    if ($cur -eq $TargetEmu) { return }

    for ($i = 1; $i -le 20000; $i++) {
        foreach ($sign in @(1, -1)) {
            $cand = $base + ($sign * $i * 0.000001)
            $Shape.Height = $cand
            $cur = [int64]([math]::Round($Shape.Height * 12700.0))
            if ($cur -eq $TargetEmu) { return }
        }
    }
    Write-Host "WARNING: could not land Height on exact EMU $TargetEmu (got $cur)"
}

function Get-ShapeById {
    param($Slide, [int]$Id)

    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $cand = $Slide.Shapes.Item($i)
        if ($cand.Id -eq $Id) { return $cand }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. Shift every existing shape right by 1270366 EMU (new X offset taken
#    straight from the target OOXML; Y/width/height are unchanged).
# ---------------------------------------------------------------------------
$shifts = @(
    @{ Id = 4; Name = 'Picture 3'; NewLeftEmu = 2203256 },
    @{ Id = 5; Name = 'Picture 4'; NewLeftEmu = 4899880 },
    @{ Id = 6; Name = 'TextBox 5'; NewLeftEmu = 2140934 },
    @{ Id = 7; Name = 'TextBox 6'; NewLeftEmu = 4572277 },
    @{ Id = 8; Name = 'Rectangle 7'; NewLeftEmu = 3744026 },
    @{ Id = 9; Name = 'Straight Arrow Connector 8'; NewLeftEmu = 3937397 },
    @{ Id = 10; Name = 'Straight Arrow Connector 9'; NewLeftEmu = 5055926 },
    @{ Id = 15; Name = 'Straight Arrow Connector 14'; NewLeftEmu = 4479502 },
    @{ Id = 16; Name = 'Straight Arrow Connector 15'; NewLeftEmu = 5417484 },
    @{ Id = 18; Name = 'Rectangle 17'; NewLeftEmu = 3434154 },
    @{ Id = 19; Name = 'Rectangle 18'; NewLeftEmu = 5147604 },
    @{ Id = 21; Name = 'Rectangle 20'; NewLeftEmu = 2684547 },
    @{ Id = 22; Name = 'Straight Arrow Connector 21'; NewLeftEmu = 4012390 },
    @{ Id = 23; Name = 'Rectangle 22'; NewLeftEmu = 2648821 },
    @{ Id = 25; Name = 'Rectangle 24'; NewLeftEmu = 2824202 },
    @{ Id = 26; Name = 'Rectangle 25'; NewLeftEmu = 4486221 },
    @{ Id = 27; Name = 'Rectangle 26'; NewLeftEmu = 3771698 },
    @{ Id = 28; Name = 'Rectangle 27'; NewLeftEmu = 6497065 },
    @{ Id = 29; Name = 'TextBox 28'; NewLeftEmu = 5362943 },
    @{ Id = 31; Name = 'TextBox 30'; NewLeftEmu = 6474308 },
    @{ Id = 33; Name = 'TextBox 32'; NewLeftEmu = 5710400 },
    @{ Id = 36; Name = 'Straight Arrow Connector 35'; NewLeftEmu = 5214978 },
    @{ Id = 37; Name = 'Straight Arrow Connector 36'; NewLeftEmu = 3199098 },
    @{ Id = 40; Name = 'Straight Arrow Connector 39'; NewLeftEmu = 3948569 },
    @{ Id = 41; Name = 'Straight Arrow Connector 40'; NewLeftEmu = 4572277 },
    @{ Id = 42; Name = 'Straight Arrow Connector 41'; NewLeftEmu = 4803505 },
    @{ Id = 45; Name = 'Straight Arrow Connector 44'; NewLeftEmu = 6647358 },
    @{ Id = 46; Name = 'Rectangle 45'; NewLeftEmu = 3305154 },
    @{ Id = 47; Name = 'Straight Arrow Connector 46'; NewLeftEmu = 3434154 },
    @{ Id = 49; Name = 'Straight Arrow Connector 48'; NewLeftEmu = 4054459 },
    @{ Id = 50; Name = 'Straight Arrow Connector 49'; NewLeftEmu = 4572277 },
    @{ Id = 57; Name = 'Rectangle 56'; NewLeftEmu = 5102593 },
    @{ Id = 58; Name = 'Straight Arrow Connector 57'; NewLeftEmu = 4726997 },
    @{ Id = 63; Name = 'Straight Arrow Connector 62'; NewLeftEmu = 6237621 },
    @{ Id = 64; Name = 'TextBox 63'; NewLeftEmu = 6654489 },
    @{ Id = 65; Name = 'TextBox 64'; NewLeftEmu = 2185772 },
    @{ Id = 66; Name = 'Straight Arrow Connector 65'; NewLeftEmu = 3169669 },
    @{ Id = 69; Name = 'Connector: Elbow 68'; NewLeftEmu = 2648821 }
)

foreach ($shift in $shifts) {
    $shape = Get-ShapeById $s $shift.Id
    if ($null -eq $shape) {
        Write-Host "WARNING: shape id" $shift.Id $shift.Name "not found"
        continue
    }
    Set-PreciseLeft $shape $shift.NewLeftEmu
}

# ---------------------------------------------------------------------------
# 2. Add the new dashed-outline "frame" rectangle behind everything else.
#    It is cloned from "Rectangle 22" (id 23), which already carries the
#    exact same no-fill / dashed-line / style / empty-centered-text
#    formatting that the new frame needs, then renamed, repositioned,
#    resized and sent to the back of the z-order.
# ---------------------------------------------------------------------------
$template = Get-ShapeById $s 23
$dup = $template.Duplicate()
$frame = $dup.Item(1)
$frame.Name = "Rectangle 1"

Set-PreciseLeft $frame 1225899
Set-PreciseTop $frame 512467
Set-PreciseWidth $frame 7998488
Set-PreciseHeight $frame 5817996

$frame.ZOrder(1)   # msoSendToBack - put it first in the shape tree
